# Applies the "feat: add 2022-Q1 data" change:
#   - the old sheet 2 ("总计") is duplicated, the duplicate becomes the new
#     "总计" (totals) sheet with an extra 2022-Q1 summary row
#   - the original sheet 2 itself is renamed to "2022-Q1" and filled with
#     the new per-fund holding breakdown for the quarter
#
# Several source columns (fund code, fund size, position %, ...) are
# numeric-looking strings that must stay TEXT (e.g. "000934" keeps its
# leading zeros). Plain `Range.Value = "000934"` auto-coerces to a number,
# and pre-setting `NumberFormat = "@"` on the destination cell leaves a
# stray text-format style behind. So `Set-AsText` stages the literal in a
# throw-away, Text-formatted scratch cell, then copies only the *value*
# (PasteSpecial xlPasteValues = -4163) into the real destination, which
# carries the text type over without dragging the scratch formatting
# along; the scratch cell is wiped immediately after with `.Clear()`.

function Set-AsText($scratch, $targetRange, $value) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $value
    $scratch.Copy()
    $targetRange.PasteSpecial(-4163)
    $scratch.Clear()
}

$wb = $excel.ActiveWorkbook
$original_active = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Duplicate the existing "总计" sheet (2nd tab) *before* touching its
#    data, so the duplicate still carries the old totals table
#    (header + the single "2021-Q3" row) and the original sheetPr block.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(2)
$q1.Copy($null, $q1)

# Free up the "总计" name on the original sheet, then claim it for the
# freshly made duplicate (sheet 3).
$q1.Name = "2022-Q1"
$total = $wb.Worksheets.Item(3)
$total.Name = "总计"

$scratch1 = $q1.Range("Z100")

# ---------------------------------------------------------------------
# 2. Turn the "2022-Q1" sheet into the per-fund holding breakdown.
# ---------------------------------------------------------------------
# Extend the header style (B1:D1 already carry the bordered header style)
# across the new columns E1:H1.
$q1.Range("B1:D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Propagate the index-column style (A2) down to the two additional rows.
$q1.Range("A2").Copy()
$q1.Range("A3:A4").PasteSpecial(-4122)

$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2

# Row 2 - 000934
Set-AsText $scratch1 $q1.Range("B2") "000934"
$q1.Range("C2").Value = "国富大中华精选混合QDII"
Set-AsText $scratch1 $q1.Range("D2") "25.71"
Set-AsText $scratch1 $q1.Range("E2") "83.59"
Set-AsText $scratch1 $q1.Range("F2") "3.78"
Set-AsText $scratch1 $q1.Range("G2") "0.9718"
$q1.Range("H2").Value = 1

# Row 3 - 006370
Set-AsText $scratch1 $q1.Range("B3") "006370"
$q1.Range("C3").Value = "国富大中华精选混合QDII美元"
Set-AsText $scratch1 $q1.Range("D3") "25.71"
Set-AsText $scratch1 $q1.Range("E3") "83.59"
Set-AsText $scratch1 $q1.Range("F3") "3.78"
Set-AsText $scratch1 $q1.Range("G3") "0.9718"
$q1.Range("H3").Value = 1

# Row 4 - 457001
Set-AsText $scratch1 $q1.Range("B4") "457001"
$q1.Range("C4").Value = "国富亚洲机会股票 (QDII)"
Set-AsText $scratch1 $q1.Range("D4") "5.93"
Set-AsText $scratch1 $q1.Range("E4") "77.36"
Set-AsText $scratch1 $q1.Range("F4") "3.40"
Set-AsText $scratch1 $q1.Range("G4") "0.2016"
$q1.Range("H4").Value = 2

# ---------------------------------------------------------------------
# 3. Update the new "总计" sheet: insert the 2022-Q1 summary row at the
#    top (row 2) and push the existing 2021-Q3 row down to row 3.
# ---------------------------------------------------------------------
# Row 3 needs the same index-column style the (only) existing data row
# (row 2, A2) already carries.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.22

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 2.15

# Restore the original active sheet/selection (the original workbook had
# "2021-Q3" selected; copying/renaming sheets above shifted focus away).
$original_active.Activate()

Write-Output "edit applied"
